$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# --- Row 1: new header labels for the lower/upper-limit theory table ---
$ws.Range("H1").Value = "lower limit"
$ws.Range("I1").Value = "upper limit"

# --- Row 2: R_convwater / 1/hA ---
$ws.Range("F2").Value = "R_convwater"
$ws.Range("G2").Value = "1/hA"
$ws.Range("H2").Formula = "=1/(`$G`$15*`$B`$14*10^-4)"
$ws.Range("I2").Formula = "=1/(G16*B14*10^-4)"

# --- Row 3: lambda_water ---
$ws.Range("F3").Value = "lambda_water"
$ws.Range("H3").Formula = "=1/`$H`$2"
$ws.Range("I3").Formula = "=1/I2"

# --- Row 4: R_condglass / L/kA ---
$ws.Range("F4").Value = "R_condglass"
$ws.Range("G4").Value = "L/kA"
$ws.Range("H4").Formula = "=B6*10^-3/(G19*B14*10^-4)"
$ws.Range("I4").Formula = "=(B6*10^-3)/(G20*B14*10^-4)"

# --- Row 5: lambda_glass ---
$ws.Range("F5").Value = "lambda_glass"
$ws.Range("H5").Formula = "=1/H4"
$ws.Range("I5").Formula = "=1/I4"

# --- Row 6: wall thickness updated from 1.6 to 1.4, plus R_air / 1/hA ---
$ws.Range("B6").Value = 1.4
$ws.Range("F6").Value = "R_air"
$ws.Range("G6").Value = "1/hA"
$ws.Range("H6").Formula = "=1/(I15*B14*10^-4)"
$ws.Range("I6").Formula = "=1/(I16*B14*10^-4)"

# --- Row 7: lambda_air ---
$ws.Range("F7").Value = "lambda_air"
$ws.Range("H7").Formula = "=1/H6"
$ws.Range("I7").Formula = "=1/I6"

# --- Row 9: lambda bottle theory (sum of resistances) ---
$ws.Range("F9").Value = "lambda bottle theory"
$ws.Range("H9").Formula = "=1/(H2+H4+H6)"
$ws.Range("I9").Formula = "=1/(I2+I4+I6)"

# --- Row 12: corrected total-surface-area formula (fixes double counting) ---
$ws.Range("B12").Formula = "=B6*0.001/(B3*(2*(B5*C5+C5*D5)+B5*D5)*0.000001)"

# --- Row 15: h_water coefficient bounds (copied from hcheattransfer.com table, Tahoma font) ---
$ws.Range("F15").Value = "lower limit h_water"
$ws.Range("F15").Font.Name = "Tahoma"
$ws.Range("F15").Font.Size = 10
$ws.Range("G15").Value = 250
$ws.Range("H15").Value = "lower limit h_air"
$ws.Range("I15").Value = 10

# --- Row 16: upper limits ---
$ws.Range("F16").Value = "upper limit h_water"
$ws.Range("G16").Value = 750
$ws.Range("H16").Value = "upper limit h_air"
$ws.Range("I16").Value = 20

# --- Row 17: source link for h coefficients ---
$ws.Range("F17").Value = "http://www.hcheattransfer.com/coefficients.html"

# --- Row 19: k_plexi lower limit (E19 carries leftover Arial/gray formatting from the pasted source) ---
$ws.Range("E19").Font.Name = "Arial"
$ws.Range("E19").Font.Color = 2236962
$ws.Range("F19").Value = "lower limit k_plexi"
$ws.Range("G19").Value = 0.17

# --- Row 20: k_plexi upper limit ---
$ws.Range("F20").Value = "upper limit k_plexi"
$ws.Range("G20").Value = 0.19

# --- Row 21: source link for k_plexi ---
$ws.Range("F21").Value = "https://www.electronics-cooling.com/2001/05/the-thermal-conductivity-of-unfilled-plastics/"

# --- column widths for the new columns ---
$ws.Range("F1").ColumnWidth = 18.21875
$ws.Range("H1").ColumnWidth = 13.77734375
$ws.Range("I1").ColumnWidth = 12.6640625

# --- restore the active selection to I10, matching the saved workbook state ---
$ws.Range("I10").Select()
